$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for all data rows (2-37)
# from 45660 to 45661 (one day later).
$ws.Range("C2:C37").Value2 = 45661

# Rows 36 and 37 swap places: the "Beteckning" (A) and "Area (ha)" (G)
# values exchange between the two rows while the rest stays put.
$a36 = $ws.Range("A36").Value2
$a37 = $ws.Range("A37").Value2
$g36 = $ws.Range("G36").Value2
$g37 = $ws.Range("G37").Value2

$ws.Range("A36").Value2 = $a37
$ws.Range("A37").Value2 = $a36

$ws.Range("G36").Value2 = $g37
$ws.Range("G37").Value2 = $g36
